$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G and H both become 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-18: H becomes 1 (G remains unchanged)
for ($r = 4; $r -le 18; $r++) {
    $ws.Range("H$r").Value = 1
}
